$d = $word.ActiveDocument

# Append "e o tipo (poupança, bolsa, tesouro etc)" to the end of the last
# content paragraph, right before the final period, e.g.:
#   "...e data da simulação." -> "...e data da simulação e o tipo (poupança, bolsa, tesouro etc)."
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("e data da simulação.", $true, $false, $false, $false, $false, `
              $true, 1, $false, "e data da simulação e o tipo (poupança, bolsa, tesouro etc).", 2)
